$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row before the "SAP Inicio Sesion" header (currently row 14).
# This pushes the existing rows 14-24 down by one, so "SAP Inicio Sesion" becomes row 15.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).RowHeight = 14.25

# Insert a new row before the "TablaEmpleado" header group (currently row 20 after the
# previous insert shifted everything down by one: old row 19 "Lenguaje" is now row 20,
# and old row 20 "Tabla Empleado en SAP" header is now row 21). Insert before row 21.
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).RowHeight = 14.25

# Populate the two new rows. Values are entered in this order so the shared-string
# table is built up the same way the original author's edit produced it.
$ws.Range("B14").Value = "euunice3@gmail.com"
$ws.Range("C14").Value = "Correo al que se le envía notificación en caso de no haber correo nuevo de RH"
$ws.Range("A21").Value = "ToInicioSAP"
$ws.Range("C21").Value = "Correo al que se le envía notificación en caso de no haber podido iniciar sesión en SAP"
$ws.Range("A14").Value = "ToCorreoNuevo"
$ws.Range("B21").Value = "euunice3@gmail.com"

# Update the sheet view: scroll so A4 is the top-left visible cell and select A21.
$ws.Activate()
$ws.Range("A21").Select()
$excel.ActiveWindow.ScrollRow = 4

# Extend the used range with two additional blank formatted rows at the bottom (999, 1000).
$ws.Rows.Item(999).RowHeight = 14.25
$ws.Rows.Item(1000).RowHeight = 14.25
